{"js": "// Fix CAF template: the Jinja condition tags call `is_resiliation()` /\n// `is_denonciation()` as if they were methods, but they are plain\n// boolean attributes, so the stray `()` must be removed, e.g.\n//   {% if convention.is_resiliation() -%}   ->  {% if convention.is_resiliation -%}\n//   {% elif convention.is_denonciation() -%} ->  {% elif convention.is_denonciation -%}\nconst body = context.document.body;\n\nconst pairs = [\n  [\"is_resiliation()\", \"is_resiliation\"],\n  [\"is_denonciation()\", \"is_denonciation\"],\n];\n\nfor (const [needle, replacement] of pairs) {\n  const results = body.search(needle, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Fix CAF template: the Jinja condition tags call `is_resiliation()` /\n# `is_denonciation()` as if they were methods, but they are plain\n# boolean attributes, so the stray `()` must be removed, e.g.\n#   {% if convention.is_resiliation() -%}   ->  {% if convention.is_resiliation -%}\n#   {% elif convention.is_denonciation() -%} ->  {% elif convention.is_denonciation -%}\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ old = \"is_resiliation()\"; new = \"is_resiliation\" },\n    @{ old = \"is_denonciation()\"; new = \"is_denonciation\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $r.old\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 0\n    $found = $rng.Find.Execute()\n    if ($found) {\n        $rng.Text = $r.new\n    }\n}\n"}
